# Feature: add arrows (arrow_n). Fixed bugs, removed unnecessary code.
#
# The "meta" sheet (key/value pairs describing the chart) gets a new
# "style" / "default" row, inserted where the previously-empty trailing
# placeholder row (A7) used to be; a fresh empty placeholder row is
# re-created one row further down.

$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("meta")

# Fill the (previously empty) placeholder row with the new key/value pair.
$meta.Range("A7").Value = "style"
$meta.Range("B7").Value = "default"

# Re-create the empty placeholder row below it, copying the key-column
# formatting (bold orange font, style index 1) from an existing key cell
# so no new style gets minted.
$meta.Range("A6").Copy()
$meta.Range("A8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
